$d = $word.ActiveDocument

# Locate the existing ".gitignore" paragraph - new content is inserted right before it.
$gitignorePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq ".gitignore – a list of files for git to not track") {
        $gitignorePara = $p
        break
    }
}

if ($gitignorePara -eq $null) {
    throw "Could not find the .gitignore paragraph"
}

$insertPoint = $d.Range($gitignorePara.Range.Start, $gitignorePara.Range.Start)

$xmlSnippet = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r><w:rPr><w:b/></w:rPr><w:t>.</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:b/></w:rPr><w:t>gitignore</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> – a list of files for </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>git</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> to not track</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Running the smoke test</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">You should run the smoke test with no arguments. It will then output the debug information of each command parsed. If there is a </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>seg</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> fault, that probably is not good and means there was an issue parsing. Look at the command being parsed in the debug statements to fix the issue. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">Whether or not things were successful, the output file should have each command tested with its command code and arguments. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">At the end of a successful run, the last line of the output file should contain a count of the failed assertions. If an assertion failed, there will be the text “—assertion failed” in the command that failed. </w:t></w:r>
</w:p>
<w:p/>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$insertPoint.InsertXML($xmlSnippet)

# InsertXML leaves one spurious trailing empty paragraph behind (the final pasted
# paragraph mark merges with whatever used to follow the insertion point, so an
# extra <w:p/> is needed in the snippet above to keep the old ".gitignore"
# paragraph separate - remove that helper paragraph now that the split is done).
$closingPara = $null
$closingParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "At the end of a successful run*") {
        $closingPara = $p
        $closingParaIndex = $i
        break
    }
}

if ($closingPara -eq $null) {
    throw "Could not find the new closing paragraph"
}

$spurious = $d.Paragraphs($closingParaIndex + 1)
$spurious.Range.Delete()

# Move the "_GoBack" bookmark from the README.md paragraph to the end of the
# newly-inserted "At the end of a successful run..." paragraph. A zero-length
# range right at a paragraph boundary is not accepted reliably, so anchor the
# bookmark on the paragraph's final (trailing-space) character instead.
$closingPara = $d.Paragraphs($closingParaIndex)
$closingEnd = $closingPara.Range.End
$bmRange = $d.Range($closingEnd - 2, $closingEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
